$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.276.49'
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").Value = '3.105.67'
$ws.Range("E3").Value = '  +5.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.85'
$ws.Range("E5").Value = '  +2.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.11'
$ws.Range("E6").Value = '  +5.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.104.05'
$ws.Range("E8").Value = '  +5.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  +1.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("E12").Value = '  +6.51%  '
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.07'
$ws.Range("E14").Value = '  +8.53%  '
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '3.623.89'
$ws.Range("E16").Value = '  +5.19%  '
$ws.Range("D17").Value = '67.294.33'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.23'
$ws.Range("E18").Value = '  +4.21%  '
$ws.Range("D19").Value = '3.104.15'
$ws.Range("E19").Value = '  +5.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.12'
$ws.Range("E20").Value = '  +17.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '471.34'
$ws.Range("E21").Value = '  +5.95%  '
$ws.Range("E22").Value = '  +6.25%  '
$ws.Range("E23").Value = '  +4.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.76'
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("E25").Value = '  +5.70%  '
$ws.Range("E26").Value = '  +5.85%  '
$ws.Range("E27").Value = '  +2.59%  '
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("E30").Value = '  +4.53%  '
$ws.Range("E32").Value = '  +5.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.64'
$ws.Range("E33").Value = '  +5.49%  '
$ws.Range("E34").Value = '  +5.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  +3.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.93'
$ws.Range("E37").Value = '  +4.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.79'
$ws.Range("E38").Value = '  +10.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.09'
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.37'
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.317'
$ws.Range("E41").Value = '  +5.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.124'
$ws.Range("E42").Value = '  +4.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.88'
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("E44").Value = '  +4.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '397.07'
$ws.Range("E45").Value = '  +4.23%  '
$ws.Range("E46").Value = '  +3.43%  '
$ws.Range("D47").Value = '2.781.75'
$ws.Range("E47").Value = '  +2.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.96'
$ws.Range("E48").Value = '  +3.50%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.85'
$ws.Range("E50").Value = '  +7.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("E51").Value = '  +5.80%  '
